$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Extend the "complementary analysis" paragraph with the new
#    sentence describing the Trinotate report-summary script.
# ---------------------------------------------------------------
$oldText = "As a complementary analysis, we have run Trinotate against a custom database with all toxins deposited in Swissprot. Aside from this change in the database for Blast searches, the program was run with the same parameters as in the main annotation process."
$addition = " For both Trinotate runs (annotation against arthropoda and toxins sequences from Swissprot), a support script packaged with Trinotate (trinotate_report_summary.pl) was used to generate result summaries, both in textual and graphical formats."
$newText = $oldText + $addition

$found = $d.Content.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)

# ---------------------------------------------------------------
# 2) Re-scope the section bookmarks so they wrap the whole section
#    (heading + body paragraphs) instead of just the heading text.
# ---------------------------------------------------------------

# Locate the key paragraphs by their heading/marker text.
function Get-ParaByText($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $needle) {
            return $p
        }
    }
    return $null
}

$pTranscriptomeHeading = Get-ParaByText "Transcriptome sequencing and assembly"
$pBioinfoHeading       = Get-ParaByText "Bioinformatics analyses"
$pTranscriptHeading    = Get-ParaByText "Transcript annotation"
$pReferencesHeading    = Get-ParaByText "REFERENCES"
$pComplementary        = Get-ParaByText $newText

# -- transcriptome-sequencing-and-assembly: heading paragraph through
#    the end of the following "FirstParagraph" body paragraph.
$endPara = $pTranscriptomeHeading.Next()
$rng = $d.Range($pTranscriptomeHeading.Range.Start, $endPara.Range.End)
$d.Bookmarks.Add("transcriptome-sequencing-and-assembly", $rng)

# -- transcript-annotation: sub-heading paragraph through the end of
#    the closing paragraph (nested inside bioinformatics-analyses, so
#    add it first -- its bookmarkEnd must land before the outer one).
$rng = $d.Range($pTranscriptHeading.Range.Start, $pComplementary.Range.End)
$d.Bookmarks.Add("transcript-annotation", $rng)

# -- bioinformatics-analyses: heading paragraph through the end of the
#    last paragraph of that whole section (right before REFERENCES).
$rng = $d.Range($pBioinfoHeading.Range.Start, $pComplementary.Range.End)
$d.Bookmarks.Add("bioinformatics-analyses", $rng)

# -- references: REFERENCES heading through the very end of the document.
$rng = $d.Range($pReferencesHeading.Range.Start, $d.Content.End)
$d.Bookmarks.Add("references", $rng)

Write-Output "done"
